$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add two new test-data rows beneath the existing project row.
# Columns: A=projectId, B=name, C=neighbourhood, D=openingDate,
#          E=closingDate, F=officerSlots, G=visibility

# Row 3: TestProjectOne / Yishun
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "TestProjectOne"
$ws.Range("C3").Value = "Yishun"
$ws.Range("D3").Value = 45741
$ws.Range("D3").NumberFormat = "dd/MM/yyyy"
$ws.Range("E3").Value = 45773
$ws.Range("E3").NumberFormat = "dd/MM/yyyy"
$ws.Range("F3").Value = 5
$ws.Range("G3").Value = $true

# Row 4: TestProject2 / YishunAgain
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "TestProject2"
$ws.Range("C4").Value = "YishunAgain"
$ws.Range("D4").Value = 45713
$ws.Range("D4").NumberFormat = "dd/MM/yyyy"
$ws.Range("E4").Value = 45940
$ws.Range("E4").NumberFormat = "dd/MM/yyyy"
$ws.Range("F4").Value = 5
$ws.Range("G4").Value = $false
